$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# ---- LP1912: update modified cells ----
$ws1.Cells.Item(2,1).Value = "Última actualización: 11:44:49"
$ws1.Cells.Item(3,1).Value = "Total filas: 199"
$ws1.Cells.Item(43,1).Value = "06:24:16"
$ws1.Cells.Item(43,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(43,4).Value = 52
$ws1.Cells.Item(44,1).Value = "06:53:31"
$ws1.Cells.Item(44,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(44,4).Value = 23
$ws1.Cells.Item(47,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(48,3).Value = "10_OLMOS"
$ws1.Cells.Item(49,3).Value = "225_GOMEZ"
$ws1.Cells.Item(54,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(55,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(56,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(78,1).Value = "07:18:07"
$ws1.Cells.Item(78,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(78,4).Value = 65
$ws1.Cells.Item(79,1).Value = "08:16:28"
$ws1.Cells.Item(79,3).Value = "215B_EL PATO"
$ws1.Cells.Item(79,4).Value = 7
$ws1.Cells.Item(89,3).Value = "10_OLMOS"
$ws1.Cells.Item(90,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(108,1).Value = "08:52:20"
$ws1.Cells.Item(108,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(108,4).Value = 31
$ws1.Cells.Item(109,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(110,1).Value = "08:38:27"
$ws1.Cells.Item(110,3).Value = "17_ROMERO"
$ws1.Cells.Item(110,4).Value = 45
$ws1.Cells.Item(116,1).Value = "09:27:56"
$ws1.Cells.Item(116,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(116,4).Value = 7
$ws1.Cells.Item(117,1).Value = "08:38:27"
$ws1.Cells.Item(117,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(117,4).Value = 56
$ws1.Cells.Item(118,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(119,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(129,3).Value = "15_ABASTO"
$ws1.Cells.Item(130,3).Value = "10_OLMOS"
$ws1.Cells.Item(148,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(149,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(150,3).Value = "215C_EL PATO"
$ws1.Cells.Item(160,1).Value = "09:27:56"
$ws1.Cells.Item(160,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(160,4).Value = 113
$ws1.Cells.Item(161,1).Value = "10:27:08"
$ws1.Cells.Item(161,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(161,4).Value = 53
$ws1.Cells.Item(174,1).Value = "11:44:49"
$ws1.Cells.Item(174,4).Value = 7
$ws1.Cells.Item(176,1).Value = "11:44:49"
$ws1.Cells.Item(176,2).Value = "11:56"
$ws1.Cells.Item(176,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(176,4).Value = 12
$ws1.Cells.Item(177,1).Value = "11:44:49"
$ws1.Cells.Item(177,2).Value = "11:58"
$ws1.Cells.Item(177,3).Value = "225_GOMEZ"
$ws1.Cells.Item(177,4).Value = 14
$ws1.Cells.Item(178,2).Value = "11:59"
$ws1.Cells.Item(178,3).Value = "225_GOMEZ"
$ws1.Cells.Item(178,4).Value = 58
$ws1.Cells.Item(179,1).Value = "11:44:49"
$ws1.Cells.Item(179,2).Value = "12:02"
$ws1.Cells.Item(179,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(179,4).Value = 18
$ws1.Cells.Item(180,1).Value = "11:01:43"
$ws1.Cells.Item(180,2).Value = "12:06"
$ws1.Cells.Item(180,3).Value = "14_ABASTO"
$ws1.Cells.Item(180,4).Value = 65
$ws1.Cells.Item(181,1).Value = "11:44:49"
$ws1.Cells.Item(181,2).Value = "12:06"
$ws1.Cells.Item(181,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(181,4).Value = 22
$ws1.Cells.Item(182,2).Value = "12:07"
$ws1.Cells.Item(182,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(182,4).Value = 100
$ws1.Cells.Item(183,1).Value = "11:44:49"
$ws1.Cells.Item(183,2).Value = "12:08"
$ws1.Cells.Item(183,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(183,4).Value = 24
$ws1.Cells.Item(184,1).Value = "11:44:49"
$ws1.Cells.Item(184,2).Value = "12:13"
$ws1.Cells.Item(184,3).Value = "10_OLMOS"
$ws1.Cells.Item(184,4).Value = 29
$ws1.Cells.Item(185,1).Value = "11:44:49"
$ws1.Cells.Item(185,2).Value = "12:14"
$ws1.Cells.Item(185,3).Value = "17_ROMERO"
$ws1.Cells.Item(185,4).Value = 30
$ws1.Cells.Item(186,1).Value = "11:44:49"
$ws1.Cells.Item(186,2).Value = "12:16"
$ws1.Cells.Item(186,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(186,4).Value = 32
$ws1.Cells.Item(187,1).Value = "10:27:08"
$ws1.Cells.Item(187,2).Value = "12:19"
$ws1.Cells.Item(187,3).Value = "14_ABASTO"
$ws1.Cells.Item(187,4).Value = 112
$ws1.Cells.Item(188,1).Value = "11:44:49"
$ws1.Cells.Item(188,2).Value = "12:20"
$ws1.Cells.Item(188,3).Value = "14_ABASTO"
$ws1.Cells.Item(188,4).Value = 36
$ws1.Cells.Item(189,1).Value = "11:44:49"
$ws1.Cells.Item(189,2).Value = "12:20"
$ws1.Cells.Item(189,3).Value = "215A_EL PATO"
$ws1.Cells.Item(189,4).Value = 36

# ---- LP1912: append new rows ----
$ws1.Cells.Item(190,1).Value = "11:44:49"
$ws1.Cells.Item(190,2).Value = "12:21"
$ws1.Cells.Item(190,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(190,4).Value = 37
$ws1.Cells.Item(190,5).Value = "LP1912"
$ws1.Cells.Item(191,1).Value = "10:27:08"
$ws1.Cells.Item(191,2).Value = "12:21"
$ws1.Cells.Item(191,3).Value = "215A_EL PATO"
$ws1.Cells.Item(191,4).Value = 114
$ws1.Cells.Item(191,5).Value = "LP1912"
$ws1.Cells.Item(192,1).Value = "11:44:49"
$ws1.Cells.Item(192,2).Value = "12:34"
$ws1.Cells.Item(192,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(192,4).Value = 50
$ws1.Cells.Item(192,5).Value = "LP1912"
$ws1.Cells.Item(193,1).Value = "11:44:49"
$ws1.Cells.Item(193,2).Value = "12:36"
$ws1.Cells.Item(193,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(193,4).Value = 52
$ws1.Cells.Item(193,5).Value = "LP1912"
$ws1.Cells.Item(194,1).Value = "11:44:49"
$ws1.Cells.Item(194,2).Value = "12:38"
$ws1.Cells.Item(194,3).Value = "17_179 Y 38"
$ws1.Cells.Item(194,4).Value = 54
$ws1.Cells.Item(194,5).Value = "LP1912"
$ws1.Cells.Item(195,1).Value = "11:44:49"
$ws1.Cells.Item(195,2).Value = "12:39"
$ws1.Cells.Item(195,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(195,4).Value = 55
$ws1.Cells.Item(195,5).Value = "LP1912"
$ws1.Cells.Item(196,1).Value = "11:44:49"
$ws1.Cells.Item(196,2).Value = "12:41"
$ws1.Cells.Item(196,3).Value = "10_OLMOS"
$ws1.Cells.Item(196,4).Value = 57
$ws1.Cells.Item(196,5).Value = "LP1912"
$ws1.Cells.Item(197,1).Value = "11:44:49"
$ws1.Cells.Item(197,2).Value = "12:48"
$ws1.Cells.Item(197,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(197,4).Value = 64
$ws1.Cells.Item(197,5).Value = "LP1912"
$ws1.Cells.Item(198,1).Value = "11:44:49"
$ws1.Cells.Item(198,2).Value = "13:02"
$ws1.Cells.Item(198,3).Value = "15_ABASTO"
$ws1.Cells.Item(198,4).Value = 78
$ws1.Cells.Item(198,5).Value = "LP1912"
$ws1.Cells.Item(199,1).Value = "11:44:49"
$ws1.Cells.Item(199,2).Value = "13:06"
$ws1.Cells.Item(199,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(199,4).Value = 82
$ws1.Cells.Item(199,5).Value = "LP1912"
$ws1.Cells.Item(200,1).Value = "11:44:49"
$ws1.Cells.Item(200,2).Value = "13:13"
$ws1.Cells.Item(200,3).Value = "215D_EL PATO"
$ws1.Cells.Item(200,4).Value = 89
$ws1.Cells.Item(200,5).Value = "LP1912"
$ws1.Cells.Item(201,1).Value = "11:44:49"
$ws1.Cells.Item(201,2).Value = "13:19"
$ws1.Cells.Item(201,3).Value = "10_OLMOS"
$ws1.Cells.Item(201,4).Value = 95
$ws1.Cells.Item(201,5).Value = "LP1912"
$ws1.Cells.Item(202,1).Value = "11:44:49"
$ws1.Cells.Item(202,2).Value = "13:21"
$ws1.Cells.Item(202,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(202,4).Value = 97
$ws1.Cells.Item(202,5).Value = "LP1912"
$ws1.Cells.Item(203,1).Value = "11:44:49"
$ws1.Cells.Item(203,2).Value = "13:26"
$ws1.Cells.Item(203,3).Value = "14_ABASTO"
$ws1.Cells.Item(203,4).Value = 102
$ws1.Cells.Item(203,5).Value = "LP1912"
$ws1.Cells.Item(204,1).Value = "11:44:49"
$ws1.Cells.Item(204,2).Value = "13:26"
$ws1.Cells.Item(204,3).Value = "15_ABASTO"
$ws1.Cells.Item(204,4).Value = 102
$ws1.Cells.Item(204,5).Value = "LP1912"

# ---- LP1912-215: update modified cells ----
$ws2.Cells.Item(2,1).Value = "Última actualización: 11:44:49"
$ws2.Cells.Item(3,1).Value = "Total filas: 21"
$ws2.Cells.Item(23,1).Value = "11:44:49"
$ws2.Cells.Item(23,4).Value = 7
$ws2.Cells.Item(24,1).Value = "11:44:49"
$ws2.Cells.Item(24,4).Value = 36

# ---- LP1912-215: append new rows ----
$ws2.Cells.Item(26,1).Value = "11:44:49"
$ws2.Cells.Item(26,2).Value = "13:13"
$ws2.Cells.Item(26,3).Value = "215D_EL PATO"
$ws2.Cells.Item(26,4).Value = 89
$ws2.Cells.Item(26,5).Value = "LP1912"

# ---- 6203-6173: update modified cells ----
$ws3.Cells.Item(2,1).Value = "Última actualización: 11:44:49"
$ws3.Cells.Item(3,1).Value = "Total filas: 25"
$ws3.Cells.Item(28,1).Value = "11:44:49"
$ws3.Cells.Item(28,4).Value = 20
$ws3.Cells.Item(29,1).Value = "11:44:49"
$ws3.Cells.Item(29,4).Value = 69

# ---- 6203-6173: append new rows ----
$ws3.Cells.Item(30,1).Value = "11:44:49"
$ws3.Cells.Item(30,2).Value = "13:30"
$ws3.Cells.Item(30,3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(30,4).Value = 106
$ws3.Cells.Item(30,5).Value = "L6173"
